$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(27, 1).Value = "'1306"
$ws.Cells.Item(27, 2).Value = "Refacciones obsoletas, Reparación no costeable, Material no disponible en el taller, Tarjeta dañada en su totalidad, Espera de refacciones, "

$ws.Cells.Item(28, 1).Value = "'1306"
$ws.Cells.Item(28, 2).Value = "Refacciones obsoletas, Reparación no costeable, Material no disponible en el taller, Tarjeta dañada en su totalidad, Espera de refacciones, componente descontinuado, "
